# g2-moderators.xlsx: replace the small sample table with the full list of
# g2/g3 moderator variable names, drop the unused "target"/c1/c2 columns,
# and turn the header row into an AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet (this also drops column C entirely).
$ws.Cells.Clear()

# Final column A / column B values, in on-sheet row order.
$values = @(
  @("generation", "name_construct"),
  @("g2", "par"),
  @("g2", "age"),
  @("g2", "rst"),
  @("g2", "car"),
  @("g2", "etp"),
  @("g2", "res"),
  @("g2", "req"),
  @("g2", "vinc"),
  @("g2", "occ"),
  @("g2", "edu"),
  @("g2", "ses"),
  @("g2", "alc"),
  @("g2", "dru"),
  @("g2", "age"),
  @("g2", "inv"),
  @("g2", "sub"),
  @("g2", "nch"),
  @("g3", "bir"),
  @("g3", "age"),
  @("g3", "sex"),
  @("g3", "gen")
)

# "sex" and "gen" (rows 21-22) were typed before the rest of the new
# column-B values during authoring, so register those shared strings
# first to match the saved string order.
$ws.Cells.Item(21, 2).Value = "sex"
$ws.Cells.Item(22, 2).Value = "gen"

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i][0]
  $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Turn the header row into an AutoFilter and register the corresponding
# hidden, sheet-scoped _FilterDatabase defined name.
$ws.Range("A1:B1").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$B`$1")
$filterName.Visible = $false

# Leave the selection where the author left it: the first empty row
# below the data.
$ws.Range("B23").Select() | Out-Null

$wb.Save()
